$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = "r775"
$ws.Range("B12").Value = "lucky"
$ws.Range("C12").Value = "im feeling"
$ws.Range("D12").Value = "2025-10-01 14:47:40"
